$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 232.3077
$ws.Range("I11").Value = 232.3077
$ws.Range("K11").Value = 232.3077
$ws.Range("M11").Value = -92.30770000000001

$ws.Range("H94").Value = 2512.9412
$ws.Range("I94").Value = 2314.6667
$ws.Range("K94").Value = 2314.6667
$ws.Range("M94").Value = -1863.6667

$ws.Range("H138").Value = 1772.4762
$ws.Range("I138").Value = 1503.0476
$ws.Range("J138").Value = 2041.9048
$ws.Range("K138").Value = 4509.142800000001
$ws.Range("L138").Value = 6125.7144
$ws.Range("M138").Value = 630.8571999999995
$ws.Range("N138").Value = -16405.7144

$ws.Range("H141").Value = 2147.3684
$ws.Range("I141").Value = 1424.1666
$ws.Range("K141").Value = 4272.4998
$ws.Range("M141").Value = 907.5002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 29800
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 29800
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H92").Value = 39083.332
$ws.Range("J92").Value = 39083.332
$ws.Range("L92").Value = 39083.332
$ws.Range("N92").Value = -44075.332

$ws.Range("H132").Value = 4257.854
$ws.Range("I132").Value = 5163.115
$ws.Range("J132").Value = 3188
$ws.Range("K132").Value = 15489.345
$ws.Range("L132").Value = 9564
$ws.Range("M132").Value = -12959.345
$ws.Range("N132").Value = -14624

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 570.26086
$ws.Range("I94").Value = 567.4286
$ws.Range("J94").Value = 600
$ws.Range("K94").Value = 567.4286
$ws.Range("L94").Value = 600
$ws.Range("M94").Value = -116.4286
$ws.Range("N94").Value = -1502

$ws.Range("H135").Value = 49703.332
$ws.Range("J135").Value = 49703.332
$ws.Range("L135").Value = 49703.332
$ws.Range("N135").Value = -59843.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 31084.285
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 31084.285
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 31084.285
$ws.Range("N74").Value = -32832.285
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 31084.285
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 31084.285
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 93252.855
$ws.Range("N77").Value = -101988.855
$ws.Range("M77").ClearContents()

$ws.Range("H92").Value = 48500
$ws.Range("J92").Value = 48500
$ws.Range("L92").Value = 48500
$ws.Range("N92").Value = -53492

$ws.Range("H132").Value = 2737.0645
$ws.Range("I132").Value = 2220.5908
$ws.Range("J132").Value = 3999.5557
$ws.Range("K132").Value = 6661.7724
$ws.Range("L132").Value = 11998.6671
$ws.Range("M132").Value = -4131.7724
$ws.Range("N132").Value = -17058.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1298.9764
$ws.Range("I68").Value = 1237.5116
$ws.Range("J68").Value = 1361.9048
$ws.Range("K68").Value = 3712.5348
$ws.Range("L68").Value = 4085.7144
$ws.Range("M68").Value = -2901.5348
$ws.Range("N68").Value = -5707.7144

$ws.Range("H71").Value = 1298.9764
$ws.Range("I71").Value = 1237.5116
$ws.Range("J71").Value = 1361.9048
$ws.Range("K71").Value = 11137.6044
$ws.Range("L71").Value = 12257.1432
$ws.Range("M71").Value = -7081.6044
$ws.Range("N71").Value = -20369.1432

$ws.Range("H97").Value = 325.36365
$ws.Range("I97").Value = 268.42856
$ws.Range("J97").Value = 425
$ws.Range("K97").Value = 805.28568
$ws.Range("L97").Value = 1275
$ws.Range("M97").Value = -309.28568
$ws.Range("N97").Value = -2267

$ws.Range("H107").Value = 996.625
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 1042
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 3126
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -6966

$ws.Range("H122").Value = 964.04877
$ws.Range("I122").Value = 808.8461
$ws.Range("J122").Value = 1233.0667
$ws.Range("K122").Value = 7279.6149
$ws.Range("L122").Value = 11097.6003
$ws.Range("M122").Value = -4829.6149
$ws.Range("N122").Value = -15997.6003

$ws.Range("H131").Value = 6657633
$ws.Range("I131").Value = 12622.223
$ws.Range("J131").Value = 11641391
$ws.Range("K131").Value = 37866.669
$ws.Range("L131").Value = 34924173
$ws.Range("M131").Value = -32826.669
$ws.Range("N131").Value = -34934253

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H122").Value = 8922.111000000001
$ws.Range("I122").Value = 56300
$ws.Range("J122").Value = 2999.875
$ws.Range("K122").Value = 168900
$ws.Range("L122").Value = 8999.625
$ws.Range("M122").Value = -166450
$ws.Range("N122").Value = -13899.625

$ws.Range("H132").Value = 65436.656
$ws.Range("I132").Value = 82455.08
$ws.Range("J132").Value = 4656.5713
$ws.Range("K132").Value = 247365.24
$ws.Range("L132").Value = 13969.7139
$ws.Range("M132").Value = -244835.24
$ws.Range("N132").Value = -19029.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 478
$ws.Range("I82").Value = 478
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 478
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -117
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 478
$ws.Range("I85").Value = 478
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 478
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 770
$ws.Range("N85").ClearContents()

$ws.Range("H93").Value = 1690714.8
$ws.Range("I93").Value = 3004614.5
$ws.Range("J93").Value = 1414.8572
$ws.Range("K93").Value = 3004614.5
$ws.Range("L93").Value = 1414.8572
$ws.Range("M93").Value = -3003366.5
$ws.Range("N93").Value = -3910.8572

$ws.Range("H132").Value = 10065.08
$ws.Range("I132").Value = 11481.45
$ws.Range("J132").Value = 4399.6
$ws.Range("K132").Value = 34444.35000000001
$ws.Range("L132").Value = 13198.8
$ws.Range("M132").Value = -31914.35000000001
$ws.Range("N132").Value = -18258.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1585.3143
$ws.Range("I132").Value = 1255
$ws.Range("J132").Value = 2700.125
$ws.Range("K132").Value = 3765
$ws.Range("L132").Value = 8100.375
$ws.Range("M132").Value = -1235
$ws.Range("N132").Value = -13160.375

$ws.Range("H136").Value = 3442.3572
$ws.Range("I136").Value = 3974.3
$ws.Range("J136").Value = 2112.5
$ws.Range("K136").Value = 11922.9
$ws.Range("L136").Value = 6337.5
$ws.Range("M136").Value = -9372.900000000001
$ws.Range("N136").Value = -11437.5
